$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Renumber the "Else, no interaction" step from "5.3)" to "5.2)".
#    The run containing ".3) " is split into ".", "2" and ") " so the "2"
#    stands on its own, matching how the digit was retyped.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Else, no interaction", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$targetCell = $findRange.Cells.Item(1)
$targetParagraph = $targetCell.Range.Paragraphs.Item(1)

$newParagraphXml = '<w:p ' + `
    'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
    'w14:paraId="06AF660E" w14:textId="13C9721A" w:rsidR="00FC21D0" ' + `
    'w:rsidRDefault="007E55D3" w:rsidP="00156273">' + `
        '<w:r><w:t>5</w:t></w:r>' + `
        '<w:r><w:t>.</w:t></w:r>' + `
        '<w:r><w:t>2</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">) </w:t></w:r>' + `
        '<w:r><w:t>Else, no interaction</w:t></w:r>' + `
        '<w:r><w:t>.</w:t></w:r>' + `
    '</w:p>'

$targetParagraph.Range.InsertXML($newParagraphXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the last row of the table ("Repeat steps 1 to 5 for all articles
#    read."), which is no longer part of the typical course of events.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Item($table.Rows.Count)
$lastRow.Delete()
